$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 (pushes existing row 4 -> row 5, etc.)
$ws.Rows("4:4").Insert()

# --- New data rows (trial 4 = lasso, trial 5 = lasso_drop, trial 6 = blend 4,5) ---
# Populate in the exact order needed so new shared strings land at the right indices:
# 11=lasso_drop, 12=lasso, 13=Henceforth..., 14=blend 4,5
$ws.Range("B7").Value = "lasso_drop"
$ws.Range("B6").Value = "lasso"
$ws.Range("B4").Value = "Henceforth, all have sigmoids"
$ws.Range("B8").Value = "blend 4,5"

# Row 6 - trial 4 - lasso
$ws.Range("A6").Value = 4
$ws.Range("C6").Value = 0.050599999999999999
$ws.Range("D6").Value = -0.36285563865669301

# Row 7 - trial 5 - lasso_drop
$ws.Range("A7").Value = 5
$ws.Range("C7").Value = 0.058599999999999999
$ws.Range("D7").Value = -2.06886488756439

# Row 8 - trial 6 - blend 4,5
$ws.Range("A8").Value = 6
$ws.Range("D8").Value = 0.046668160474799401

# Style the new note row (B4): italic, red font
$ws.Range("B4").Font.Italic = $true
$ws.Range("B4").Font.Color = 255

# Update the active selection to D9, matching the post-edit workbook state
$ws.Range("D9").Select()
